# Add a "Source:" attribution textbox under the Azure graphic on slide 6,
# and move the picture up so the offset y becomes 0.
#
# NOTE: PowerPoint's COM object model expresses Shape geometry (Left/Top/
# Width/Height, and AddTextbox's position/size args) in POINTS, while the
# underlying OOXML <a:off>/<a:ext> store EMU (1 pt = 12700 EMU). To land on
# an exact EMU value we divide by 12700 before handing it to the COM layer.
$EMU_PER_PT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# --- Move the picture's vertical offset to 0 (was 213081 EMU) -------------
$pic = $s.Shapes.Item(1)
$pic.Top = 0 / $EMU_PER_PT

# --- Add the new "Source: ..." textbox -------------------------------------
$left   = 221673    / $EMU_PER_PT
$top    = 6483925   / $EMU_PER_PT
$width  = 10871200  / $EMU_PER_PT
$height = 369332    / $EMU_PER_PT

$textBox = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)

$textBox.TextFrame.WordWrap = -1
$textBox.TextFrame.AutoSize = 1

$url = "http://msdn.microsoft.com/en-us/library/azure/dd163896.aspx"

$tr = $textBox.TextFrame.TextRange
$tr.Text = "Source: http://msdn.microsoft.com/en-us/library/azure/dd163896.aspx "

# Split out the hyperlinked runs exactly like PowerPoint does when it
# autolinks a typed URL: "http" / "://" / "msdn.microsoft.com/...aspx".
$httpRun  = $tr.Characters(9, 4)
$httpRun.ActionSettings(1).Hyperlink.Address = $url

$slashRun = $tr.Characters(13, 3)
$slashRun.ActionSettings(1).Hyperlink.Address = $url

$restRun  = $tr.Characters(16, 52)
$restRun.ActionSettings(1).Hyperlink.Address = $url
